# Apply the Efna4-Epha2 NATMI LR-pair update: add the "M2" target/sending cluster
# and refresh all derived statistics across the 4x4 ECs/FAPs/M2/sCs grid (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "Efna4", "Epha2", "ECs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 24.244885, 72.734655, 0.6895205882382217, 0.6895205882382218, 20.21888445806833, 181.969960122615, 0.210443912702364, 0.2104439127023641),
    @("ECs", "Efna4", "Epha2", "FAPs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 0.7651789999999999, 2.295537, 0.02176156637523891, 0.02176156637523891, 0.6381166910356666, 5.743050219321, 0.00664170041135201, 0.006641700411352013),
    @("ECs", "Efna4", "Epha2", "M2", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 0.6351283333333334, 1.905385, 0.01806294655581008, 0.01806294655581008, 0.5296616745227778, 4.766955070705, 0.005512869685081945, 0.005512869685081946),
    @("ECs", "Efna4", "Epha2", "sCs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 9.516752666666667, 28.550258, 0.2706548988307292, 0.2706548988307293, 7.936441958101556, 71.427977622914, 0.08260475013158404, 0.08260475013158407),
    @("FAPs", "Efna4", "Epha2", "ECs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 24.244885, 72.734655, 0.6895205882382217, 0.6895205882382218, 37.27744007475833, 335.496960672825, 0.3879942219922812, 0.3879942219922813),
    @("FAPs", "Efna4", "Epha2", "FAPs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 0.7651789999999999, 2.295537, 0.02176156637523891, 0.02176156637523891, 1.176492044361666, 10.588428399255, 0.01224526454919591, 0.01224526454919592),
    @("FAPs", "Efna4", "Epha2", "M2", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 0.6351283333333334, 1.905385, 0.01806294655581008, 0.01806294655581008, 0.9765341590861112, 8.788807431775, 0.01016404588254063, 0.01016404588254063),
    @("FAPs", "Efna4", "Epha2", "sCs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 9.516752666666667, 28.550258, 0.2706548988307292, 0.2706548988307293, 14.63237203385222, 131.69134830467, 0.1522978989917379, 0.152297898991738),
    @("M2", "Efna4", "Epha2", "ECs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 24.244885, 72.734655, 0.6895205882382217, 0.6895205882382218, 0.9671688707916667, 8.704519837125, 0.01006656929245694, 0.01006656929245694),
    @("M2", "Efna4", "Epha2", "FAPs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 0.7651789999999999, 2.295537, 0.02176156637523891, 0.02176156637523891, 0.03052426560833333, 0.274718390475, 0.0003177052571968441, 0.0003177052571968442),
    @("M2", "Efna4", "Epha2", "M2", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 0.6351283333333334, 1.905385, 0.01806294655581008, 0.01806294655581008, 0.02533632776388889, 0.228026949875, 0.0002637077213235983, 0.0002637077213235984),
    @("M2", "Efna4", "Epha2", "sCs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 9.516752666666667, 28.550258, 0.2706548988307292, 0.2706548988307293, 0.3796391251277778, 3.41675212615, 0.003951392228017346, 0.003951392228017348),
    @("sCs", "Efna4", "Epha2", "ECs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 24.244885, 72.734655, 0.6895205882382217, 0.6895205882382218, 7.783788002736667, 70.05409202463001, 0.08101588425111954, 0.08101588425111957),
    @("sCs", "Efna4", "Epha2", "FAPs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 0.7651789999999999, 2.295537, 0.02176156637523891, 0.02176156637523891, 0.2456596977113333, 2.210937279402, 0.002556896157494143, 0.002556896157494144),
    @("sCs", "Efna4", "Epha2", "M2", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 0.6351283333333334, 1.905385, 0.01806294655581008, 0.01806294655581008, 0.2039071045788889, 1.83516394121, 0.00212232326686391, 0.00212232326686391),
    @("sCs", "Efna4", "Epha2", "sCs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 9.516752666666667, 28.550258, 0.2706548988307292, 0.2706548988307293, 3.055340754629778, 27.498066791668, 0.03180085747938997, 0.03180085747938999)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $rows[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}

Write-Host "Updated sheet1 rows 2-17 (A1:T17)"